$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.850.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.14%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.815.91'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.16%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.48%  '

# Row 5
$ws.Range("E5").Value = '  -0.33%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '308.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.90%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4621'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.51%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3637'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.39%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07221'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.16%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8578'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.02%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.72'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.55%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07526'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.64%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.790.80'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -9.39%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.320'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.35%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.533'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.67%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.74'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.58%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.008'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.17%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008565'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.79%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.008'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.24%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.162.53'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.13%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.41'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.66%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.143'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.43%  '

# Row 23
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.50'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.77%  '

# Row 24
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.116.68'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.86%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.31'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.50%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.850'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.23%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.75%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.064'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.76%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.080'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.20%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '114.99'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.45%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08857'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.58%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.950'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.01%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.408'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.28%  '

# Row 34
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.130'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.33%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7158'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.51%  '

# Row 36
$ws.Range("E36").Value = '  -0.40%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.072'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.01%  '

# Row 38
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.436'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.53%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05233'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.92%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01911'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.20%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.919'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.65%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.146'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.63%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5125'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.86%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1622'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.29%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.170'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.00%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4793'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.52%  '

# Row 47
$ws.Range("E47").Value = '  -0.36%  '

# Row 48
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.97'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.95%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.04'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.18%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.614'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.83%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06198'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.91%  '
